# Apply the updated NATMI Col1a2-Cd93 LR-pair numbers (20 sending/target-cluster
# combinations, incl. the 4 new "sCs" target rows) to Sheet1, rows 2-21.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = New-Object "object[,]" 20,20
$data[0,0] = "ECs"
$data[0,1] = "Col1a2"
$data[0,2] = "Cd93"
$data[0,3] = "ECs"
$data[0,4] = [double]"3.0"
$data[0,5] = [double]"1.0"
$data[0,6] = [double]"65.399996"
$data[0,7] = [double]"196.199988"
$data[0,8] = [double]"0.0234573392998008"
$data[0,9] = [double]"0.02345733929980081"
$data[0,10] = [double]"3.0"
$data[0,11] = [double]"1.0"
$data[0,12] = [double]"229.5846506666667"
$data[0,13] = [double]"688.753952"
$data[0,14] = [double]"0.5033187360873315"
$data[0,15] = [double]"0.5033187360873315"
$data[0,16] = [double]"15014.8352352614"
$data[0,17] = [double]"135133.5171173526"
$data[0,18] = [double]"0.01180651836834743"
$data[0,19] = [double]"0.01180651836834743"

$data[1,0] = "ECs"
$data[1,1] = "Col1a2"
$data[1,2] = "Cd93"
$data[1,3] = "M1"
$data[1,4] = [double]"3.0"
$data[1,5] = [double]"1.0"
$data[1,6] = [double]"65.399996"
$data[1,7] = [double]"196.199988"
$data[1,8] = [double]"0.0234573392998008"
$data[1,9] = [double]"0.02345733929980081"
$data[1,10] = [double]"3.0"
$data[1,11] = [double]"1.0"
$data[1,12] = [double]"135.7283196666666"
$data[1,13] = [double]"407.1849589999999"
$data[1,14] = [double]"0.2975573763642838"
$data[1,15] = [double]"0.2975573763642838"
$data[1,16] = [double]"8876.631563286719"
$data[1,17] = [double]"79889.68406958049"
$data[1,18] = [double]"0.006979904338535532"
$data[1,19] = [double]"0.006979904338535534"

$data[2,0] = "ECs"
$data[2,1] = "Col1a2"
$data[2,2] = "Cd93"
$data[2,3] = "M2"
$data[2,4] = [double]"3.0"
$data[2,5] = [double]"1.0"
$data[2,6] = [double]"65.399996"
$data[2,7] = [double]"196.199988"
$data[2,8] = [double]"0.0234573392998008"
$data[2,9] = [double]"0.02345733929980081"
$data[2,10] = [double]"3.0"
$data[2,11] = [double]"1.0"
$data[2,12] = [double]"90.23148833333335"
$data[2,13] = [double]"270.694465"
$data[2,14] = [double]"0.1978146123067711"
$data[2,15] = [double]"0.1978146123067711"
$data[2,16] = [double]"5901.138976074048"
$data[2,17] = [double]"53110.25078466643"
$data[2,18] = [double]"0.004640204479338481"
$data[2,19] = [double]"0.004640204479338481"

$data[3,0] = "ECs"
$data[3,1] = "Col1a2"
$data[3,2] = "Cd93"
$data[3,3] = "sCs"
$data[3,4] = [double]"3.0"
$data[3,5] = [double]"1.0"
$data[3,6] = [double]"65.399996"
$data[3,7] = [double]"196.199988"
$data[3,8] = [double]"0.0234573392998008"
$data[3,9] = [double]"0.02345733929980081"
$data[3,10] = [double]"3.0"
$data[3,11] = [double]"1.0"
$data[3,12] = [double]"0.5972149999999999"
$data[3,13] = [double]"1.791645"
$data[3,14] = [double]"0.001309275241613694"
$data[3,15] = [double]"0.001309275241613694"
$data[3,16] = [double]"39.05785861114"
$data[3,17] = [double]"351.52072750026"
$data[3,18] = [double]"3.071211357936111e-05"
$data[3,19] = [double]"3.071211357936112e-05"

$data[4,0] = "FAPs"
$data[4,1] = "Col1a2"
$data[4,2] = "Cd93"
$data[4,3] = "ECs"
$data[4,4] = [double]"3.0"
$data[4,5] = [double]"1.0"
$data[4,6] = [double]"2623.51945"
$data[4,7] = [double]"7870.55835"
$data[4,8] = [double]"0.9409906676183405"
$data[4,9] = [double]"0.9409906676183405"
$data[4,10] = [double]"3.0"
$data[4,11] = [double]"1.0"
$data[4,12] = [double]"229.5846506666667"
$data[4,13] = [double]"688.753952"
$data[4,14] = [double]"0.5033187360873315"
$data[4,15] = [double]"0.5033187360873315"
$data[4,16] = [double]"602319.7964454555"
$data[4,17] = [double]"5420878.1680091"
$data[4,18] = [double]"0.4736182334956374"
$data[4,19] = [double]"0.4736182334956374"

$data[5,0] = "FAPs"
$data[5,1] = "Col1a2"
$data[5,2] = "Cd93"
$data[5,3] = "M1"
$data[5,4] = [double]"3.0"
$data[5,5] = [double]"1.0"
$data[5,6] = [double]"2623.51945"
$data[5,7] = [double]"7870.55835"
$data[5,8] = [double]"0.9409906676183405"
$data[5,9] = [double]"0.9409906676183405"
$data[5,10] = [double]"3.0"
$data[5,11] = [double]"1.0"
$data[5,12] = [double]"135.7283196666666"
$data[5,13] = [double]"407.1849589999999"
$data[5,14] = [double]"0.2975573763642838"
$data[5,15] = [double]"0.2975573763642838"
$data[5,16] = [double]"356085.8865613175"
$data[5,17] = [double]"3204772.979051857"
$data[5,18] = [double]"0.2799987142397892"
$data[5,19] = [double]"0.2799987142397892"

$data[6,0] = "FAPs"
$data[6,1] = "Col1a2"
$data[6,2] = "Cd93"
$data[6,3] = "M2"
$data[6,4] = [double]"3.0"
$data[6,5] = [double]"1.0"
$data[6,6] = [double]"2623.51945"
$data[6,7] = [double]"7870.55835"
$data[6,8] = [double]"0.9409906676183405"
$data[6,9] = [double]"0.9409906676183405"
$data[6,10] = [double]"3.0"
$data[6,11] = [double]"1.0"
$data[6,12] = [double]"90.23148833333335"
$data[6,13] = [double]"270.694465"
$data[6,14] = [double]"0.1978146123067711"
$data[6,15] = [double]"0.1978146123067711"
$data[6,16] = [double]"236724.0646449481"
$data[6,17] = [double]"2130516.581804533"
$data[6,18] = [double]"0.1861417040992117"
$data[6,19] = [double]"0.1861417040992117"

$data[7,0] = "FAPs"
$data[7,1] = "Col1a2"
$data[7,2] = "Cd93"
$data[7,3] = "sCs"
$data[7,4] = [double]"3.0"
$data[7,5] = [double]"1.0"
$data[7,6] = [double]"2623.51945"
$data[7,7] = [double]"7870.55835"
$data[7,8] = [double]"0.9409906676183405"
$data[7,9] = [double]"0.9409906676183405"
$data[7,10] = [double]"3.0"
$data[7,11] = [double]"1.0"
$data[7,12] = [double]"0.5972149999999999"
$data[7,13] = [double]"1.791645"
$data[7,14] = [double]"0.001309275241613694"
$data[7,15] = [double]"0.001309275241613694"
$data[7,16] = [double]"1566.80516833175"
$data[7,17] = [double]"14101.24651498575"
$data[7,18] = [double]"0.001232015783702234"
$data[7,19] = [double]"0.001232015783702234"

$data[8,0] = "M1"
$data[8,1] = "Col1a2"
$data[8,2] = "Cd93"
$data[8,3] = "ECs"
$data[8,4] = [double]"3.0"
$data[8,5] = [double]"1.0"
$data[8,6] = [double]"0.220081"
$data[8,7] = [double]"0.6602429999999999"
$data[8,8] = [double]"7.89375383209421e-05"
$data[8,9] = [double]"7.893753832094211e-05"
$data[8,10] = [double]"3.0"
$data[8,11] = [double]"1.0"
$data[8,12] = [double]"229.5846506666667"
$data[8,13] = [double]"688.753952"
$data[8,14] = [double]"0.5033187360873315"
$data[8,15] = [double]"0.5033187360873315"
$data[8,16] = [double]"50.52721950337066"
$data[8,17] = [double]"454.744975530336"
$data[8,18] = [double]"3.973074201754187e-05"
$data[8,19] = [double]"3.973074201754187e-05"

$data[9,0] = "M1"
$data[9,1] = "Col1a2"
$data[9,2] = "Cd93"
$data[9,3] = "M1"
$data[9,4] = [double]"3.0"
$data[9,5] = [double]"1.0"
$data[9,6] = [double]"0.220081"
$data[9,7] = [double]"0.6602429999999999"
$data[9,8] = [double]"7.89375383209421e-05"
$data[9,9] = [double]"7.893753832094211e-05"
$data[9,10] = [double]"3.0"
$data[9,11] = [double]"1.0"
$data[9,12] = [double]"135.7283196666666"
$data[9,13] = [double]"407.1849589999999"
$data[9,14] = [double]"0.2975573763642838"
$data[9,15] = [double]"0.2975573763642838"
$data[9,16] = [double]"29.87122432055966"
$data[9,17] = [double]"268.8410188850369"
$data[9,18] = [double]"2.348844679943464e-05"
$data[9,19] = [double]"2.348844679943464e-05"

$data[10,0] = "M1"
$data[10,1] = "Col1a2"
$data[10,2] = "Cd93"
$data[10,3] = "M2"
$data[10,4] = [double]"3.0"
$data[10,5] = [double]"1.0"
$data[10,6] = [double]"0.220081"
$data[10,7] = [double]"0.6602429999999999"
$data[10,8] = [double]"7.89375383209421e-05"
$data[10,9] = [double]"7.893753832094211e-05"
$data[10,10] = [double]"3.0"
$data[10,11] = [double]"1.0"
$data[10,12] = [double]"90.23148833333335"
$data[10,13] = [double]"270.694465"
$data[10,14] = [double]"0.1978146123067711"
$data[10,15] = [double]"0.1978146123067711"
$data[10,16] = [double]"19.85823618388833"
$data[10,17] = [double]"178.724125654995"
$data[10,18] = [double]"1.561499853940804e-05"
$data[10,19] = [double]"1.561499853940805e-05"

$data[11,0] = "M1"
$data[11,1] = "Col1a2"
$data[11,2] = "Cd93"
$data[11,3] = "sCs"
$data[11,4] = [double]"3.0"
$data[11,5] = [double]"1.0"
$data[11,6] = [double]"0.220081"
$data[11,7] = [double]"0.6602429999999999"
$data[11,8] = [double]"7.89375383209421e-05"
$data[11,9] = [double]"7.893753832094211e-05"
$data[11,10] = [double]"3.0"
$data[11,11] = [double]"1.0"
$data[11,12] = [double]"0.5972149999999999"
$data[11,13] = [double]"1.791645"
$data[11,14] = [double]"0.001309275241613694"
$data[11,15] = [double]"0.001309275241613694"
$data[11,16] = [double]"0.131435674415"
$data[11,17] = [double]"1.182921069735"
$data[11,18] = [double]"1.033509645575417e-07"
$data[11,19] = [double]"1.033509645575418e-07"

$data[12,0] = "M2"
$data[12,1] = "Col1a2"
$data[12,2] = "Cd93"
$data[12,3] = "ECs"
$data[12,4] = [double]"2.0"
$data[12,5] = [double]"0.6666666666666666"
$data[12,6] = [double]"0.250186"
$data[12,7] = [double]"0.7505580000000001"
$data[12,8] = [double]"8.973544723244271e-05"
$data[12,9] = [double]"8.973544723244272e-05"
$data[12,10] = [double]"3.0"
$data[12,11] = [double]"1.0"
$data[12,12] = [double]"229.5846506666667"
$data[12,13] = [double]"688.753952"
$data[12,14] = [double]"0.5033187360873315"
$data[12,15] = [double]"0.5033187360873315"
$data[12,16] = [double]"57.43886541169067"
$data[12,17] = [double]"516.9497887052161"
$data[12,18] = [double]"4.516553188326449e-05"
$data[12,19] = [double]"4.516553188326449e-05"

$data[13,0] = "M2"
$data[13,1] = "Col1a2"
$data[13,2] = "Cd93"
$data[13,3] = "M1"
$data[13,4] = [double]"2.0"
$data[13,5] = [double]"0.6666666666666666"
$data[13,6] = [double]"0.250186"
$data[13,7] = [double]"0.7505580000000001"
$data[13,8] = [double]"8.973544723244271e-05"
$data[13,9] = [double]"8.973544723244272e-05"
$data[13,10] = [double]"3.0"
$data[13,11] = [double]"1.0"
$data[13,12] = [double]"135.7283196666666"
$data[13,13] = [double]"407.1849589999999"
$data[13,14] = [double]"0.2975573763642838"
$data[13,15] = [double]"0.2975573763642838"
$data[13,16] = [double]"33.95732538412466"
$data[13,17] = [double]"305.615928457122"
$data[13,18] = [double]"2.670144424536128e-05"
$data[13,19] = [double]"2.670144424536129e-05"

$data[14,0] = "M2"
$data[14,1] = "Col1a2"
$data[14,2] = "Cd93"
$data[14,3] = "M2"
$data[14,4] = [double]"2.0"
$data[14,5] = [double]"0.6666666666666666"
$data[14,6] = [double]"0.250186"
$data[14,7] = [double]"0.7505580000000001"
$data[14,8] = [double]"8.973544723244271e-05"
$data[14,9] = [double]"8.973544723244272e-05"
$data[14,10] = [double]"3.0"
$data[14,11] = [double]"1.0"
$data[14,12] = [double]"90.23148833333335"
$data[14,13] = [double]"270.694465"
$data[14,14] = [double]"0.1978146123067711"
$data[14,15] = [double]"0.1978146123067711"
$data[14,16] = [double]"22.57465514016334"
$data[14,17] = [double]"203.17189626147"
$data[14,18] = [double]"1.775098270446037e-05"
$data[14,19] = [double]"1.775098270446037e-05"

$data[15,0] = "M2"
$data[15,1] = "Col1a2"
$data[15,2] = "Cd93"
$data[15,3] = "sCs"
$data[15,4] = [double]"2.0"
$data[15,5] = [double]"0.6666666666666666"
$data[15,6] = [double]"0.250186"
$data[15,7] = [double]"0.7505580000000001"
$data[15,8] = [double]"8.973544723244271e-05"
$data[15,9] = [double]"8.973544723244272e-05"
$data[15,10] = [double]"3.0"
$data[15,11] = [double]"1.0"
$data[15,12] = [double]"0.5972149999999999"
$data[15,13] = [double]"1.791645"
$data[15,14] = [double]"0.001309275241613694"
$data[15,15] = [double]"0.001309275241613694"
$data[15,16] = [double]"0.14941483199"
$data[15,17] = [double]"1.34473348791"
$data[15,18] = [double]"1.174883993565694e-07"
$data[15,19] = [double]"1.174883993565694e-07"

$data[16,0] = "sCs"
$data[16,1] = "Col1a2"
$data[16,2] = "Cd93"
$data[16,3] = "ECs"
$data[16,4] = [double]"3.0"
$data[16,5] = [double]"1.0"
$data[16,6] = [double]"98.65010533333333"
$data[16,7] = [double]"295.950316"
$data[16,8] = [double]"0.03538332009630534"
$data[16,9] = [double]"0.03538332009630534"
$data[16,10] = [double]"3.0"
$data[16,11] = [double]"1.0"
$data[16,12] = [double]"229.5846506666667"
$data[16,13] = [double]"688.753952"
$data[16,14] = [double]"0.5033187360873315"
$data[16,15] = [double]"0.5033187360873315"
$data[16,16] = [double]"22648.5499711832"
$data[16,17] = [double]"203836.9497406488"
$data[16,18] = [double]"0.01780908794944588"
$data[16,19] = [double]"0.01780908794944588"

$data[17,0] = "sCs"
$data[17,1] = "Col1a2"
$data[17,2] = "Cd93"
$data[17,3] = "M1"
$data[17,4] = [double]"3.0"
$data[17,5] = [double]"1.0"
$data[17,6] = [double]"98.65010533333333"
$data[17,7] = [double]"295.950316"
$data[17,8] = [double]"0.03538332009630534"
$data[17,9] = [double]"0.03538332009630534"
$data[17,10] = [double]"3.0"
$data[17,11] = [double]"1.0"
$data[17,12] = [double]"135.7283196666666"
$data[17,13] = [double]"407.1849589999999"
$data[17,14] = [double]"0.2975573763642838"
$data[17,15] = [double]"0.2975573763642838"
$data[17,16] = [double]"13389.613031833"
$data[17,17] = [double]"120506.517286497"
$data[17,18] = [double]"0.01052856789491425"
$data[17,19] = [double]"0.01052856789491425"

$data[18,0] = "sCs"
$data[18,1] = "Col1a2"
$data[18,2] = "Cd93"
$data[18,3] = "M2"
$data[18,4] = [double]"3.0"
$data[18,5] = [double]"1.0"
$data[18,6] = [double]"98.65010533333333"
$data[18,7] = [double]"295.950316"
$data[18,8] = [double]"0.03538332009630534"
$data[18,9] = [double]"0.03538332009630534"
$data[18,10] = [double]"3.0"
$data[18,11] = [double]"1.0"
$data[18,12] = [double]"90.23148833333335"
$data[18,13] = [double]"270.694465"
$data[18,14] = [double]"0.1978146123067711"
$data[18,15] = [double]"0.1978146123067711"
$data[18,16] = [double]"8901.345828466772"
$data[18,17] = [double]"80112.11245620095"
$data[18,18] = [double]"0.006999337746977022"
$data[18,19] = [double]"0.006999337746977022"

$data[19,0] = "sCs"
$data[19,1] = "Col1a2"
$data[19,2] = "Cd93"
$data[19,3] = "sCs"
$data[19,4] = [double]"3.0"
$data[19,5] = [double]"1.0"
$data[19,6] = [double]"98.65010533333333"
$data[19,7] = [double]"295.950316"
$data[19,8] = [double]"0.03538332009630534"
$data[19,9] = [double]"0.03538332009630534"
$data[19,10] = [double]"3.0"
$data[19,11] = [double]"1.0"
$data[19,12] = [double]"0.5972149999999999"
$data[19,13] = [double]"1.791645"
$data[19,14] = [double]"0.001309275241613694"
$data[19,15] = [double]"0.001309275241613694"
$data[19,16] = [double]"58.91532265664666"
$data[19,17] = [double]"530.2379039098199"
$data[19,18] = [double]"4.632650496818486e-05"
$data[19,19] = [double]"4.632650496818486e-05"

$ws.Range("A2:T21").Value = $data

